$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'33.896.11"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Formula = "'1.778.58"
$ws.Range("E3").Value = "  -0.67%  "

$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").Formula = "'220.50"
$ws.Range("E5").Value = "  -1.35%  "

$ws.Range("D6").Formula = "'0.550"
$ws.Range("E6").Value = "  -0.52%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("D8").Formula = "'31.94"
$ws.Range("E8").Value = "  -0.87%  "

$ws.Range("D9").Formula = "'0.281"
$ws.Range("E9").Value = "  -0.17%  "

$ws.Range("D10").Formula = "'0.0706"
$ws.Range("E10").Value = "  -0.07%  "

$ws.Range("D11").Formula = "'0.0928"

$ws.Range("D12").Formula = "'2.040.23"
$ws.Range("E12").Value = "  -0.32%  "

$ws.Range("D13").Formula = "'1.788.60"
$ws.Range("E13").Value = "  -0.03%  "

$ws.Range("D14").Formula = "'10.68"
$ws.Range("E14").Value = "  -3.03%  "

$ws.Range("D15").Formula = "'0.620"
$ws.Range("E15").Value = "  -2.06%  "

$ws.Range("D16").Formula = "'33.853.01"
$ws.Range("E16").Value = "  -2.32%  "

$ws.Range("D17").Formula = "'4.13"
$ws.Range("E17").Value = "  -4.04%  "

$ws.Range("D18").Formula = "'67.46"
$ws.Range("E18").Value = "  -2.35%  "

$ws.Range("D19").Formula = "'242.16"
$ws.Range("E19").Value = "  -4.67%  "

$ws.Range("D20").Formula = "'0.0₃0775"
$ws.Range("E20").Value = "  -3.78%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").Formula = "'10.65"
$ws.Range("E22").Value = "  +0.70%  "

$ws.Range("D23").Formula = "'4.03"
$ws.Range("E23").Value = "  -3.47%  "

$ws.Range("D24").Formula = "'2.10"
$ws.Range("E24").Value = "  -1.65%  "

$ws.Range("D25").Formula = "'157.24"
$ws.Range("E25").Value = "  -1.86%  "

$ws.Range("D26").Formula = "'16.23"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("D27").Formula = "'6.97"
$ws.Range("E27").Value = "  -1.98%  "

$ws.Range("D28").Formula = "'0.112"
$ws.Range("E28").Value = "  -2.12%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").Formula = "'0.0515"
$ws.Range("E30").Value = "  -2.18%  "

$ws.Range("E31").Value = "  -0.02%  "

$ws.Range("D32").Formula = "'3.64"
$ws.Range("E32").Value = "  -3.58%  "

$ws.Range("D33").Formula = "'3.46"
$ws.Range("E33").Value = "  -4.12%  "

$ws.Range("D34").Formula = "'1.80"
$ws.Range("E34").Value = "  -3.64%  "

$ws.Range("D35").Formula = "'1.387.16"
$ws.Range("E35").Value = "  -3.34%  "

$ws.Range("D36").Formula = "'0.632"
$ws.Range("E36").Value = "  -0.73%  "

$ws.Range("D37").Formula = "'1.04"
$ws.Range("E37").Value = "  -0.97%  "

$ws.Range("D38").Formula = "'0.0184"
$ws.Range("E38").Value = "  -4.09%  "

$ws.Range("B39").Value = "HuobiToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D39").Formula = "'2.35"
$ws.Range("E39").Value = "  +1.15%  "

$ws.Range("D40").Formula = "'78.64"
$ws.Range("E40").Value = "  -7.18%  "

$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Formula = "'0.912"
$ws.Range("E41").Value = "  -1.25%  "

$ws.Range("D42").Formula = "'2.69"
$ws.Range("E42").Value = "  -3.33%  "

$ws.Range("D43").Formula = "'2.11"
$ws.Range("E43").Value = "  +1.56%  "

$ws.Range("B44").Value = "WEMIXToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D44").Formula = "'1.05"
$ws.Range("E44").Value = "  -0.72%  "

$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").Formula = "'0.0494"
$ws.Range("E45").Value = "  +1.60%  "

$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Formula = "'5.85"
$ws.Range("E46").Value = "  -1.90%  "

$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").Formula = "'106.22"
$ws.Range("E47").Value = "  +0.61%  "

$ws.Range("D48").Formula = "'1.938.09"
$ws.Range("E48").Value = "  -0.18%  "

$ws.Range("E49").Value = "  +0.14%  "

$ws.Range("D50").Formula = "'11.80"
$ws.Range("E50").Value = "  -1.09%  "

$ws.Range("E51").Value = "  +0.42%  "
